# chore: update Sheets via scheduled runner
# Refreshes cached market-price / leve-profit figures (columns H-N) on a
# handful of rows across the per-job sheets (ALC, ARM, BSM, CRP, CUL, GSM,
# LTW, WVR). A couple of rows pick up a brand-new M/N value where none was
# cached before, and CRP!N107 loses its stale cached value entirely.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 870.7692
$ws.Range("I98").Value = 692.8333
$ws.Range("K98").Value = 692.8333
$ws.Range("M98").Value = 805.1667
$ws.Range("H107").Value = 951.6667
$ws.Range("I107").Value = 951.6667
$ws.Range("K107").Value = 951.6667
$ws.Range("M107").Value = 968.3333
$ws.Range("H122").Value = 870.7692
$ws.Range("I122").Value = 692.8333
$ws.Range("K122").Value = 2078.4999
$ws.Range("M122").Value = 371.5001000000002
$ws.Range("H135").Value = 2137.3333
$ws.Range("I135").Value = 1342.6
$ws.Range("J135").Value = 3130.75
$ws.Range("K135").Value = 12083.4
$ws.Range("L135").Value = 28176.75
$ws.Range("M135").Value = -9548.4
$ws.Range("N135").Value = -33246.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5000
$ws.Range("I2").Value = 5000
$ws.Range("K2").Value = 5000
$ws.Range("M2").Value = -4887
$ws.Range("H74").Value = 1229.1482
$ws.Range("I74").Value = 961.2763
$ws.Range("K74").Value = 961.2763
$ws.Range("M74").Value = -87.27629999999999
$ws.Range("H77").Value = 1229.1482
$ws.Range("I77").Value = 961.2763
$ws.Range("K77").Value = 4806.3815
$ws.Range("M77").Value = -438.3814999999995
$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 5000
$ws.Range("K116").Value = 5000
$ws.Range("M116").Value = -2706
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 5000
$ws.Range("K3").Value = 5000
$ws.Range("M3").Value = -4886
$ws.Range("H80").Value = 618.75
$ws.Range("J80").Value = 459.5
$ws.Range("L80").Value = 459.5
$ws.Range("N80").Value = -2455.5
$ws.Range("H83").Value = 618.75
$ws.Range("J83").Value = 459.5
$ws.Range("L83").Value = 2297.5
$ws.Range("N83").Value = -12281.5
$ws.Range("H86").Value = 1933.875
$ws.Range("J86").Value = 2802
$ws.Range("L86").Value = 2802
$ws.Range("N86").Value = -5048
$ws.Range("H89").Value = 1933.875
$ws.Range("J89").Value = 2802
$ws.Range("L89").Value = 14010
$ws.Range("N89").Value = -25242
$ws.Range("H103").Value = 8977.75
$ws.Range("J103").Value = 8977.75
$ws.Range("L103").Value = 8977.75
$ws.Range("N103").Value = -11321.75
$ws.Range("H134").Value = 2467.9473
$ws.Range("I134").Value = 2126.2666
$ws.Range("J134").Value = 3749.25
$ws.Range("K134").Value = 6378.7998
$ws.Range("L134").Value = 11247.75
$ws.Range("M134").Value = -3843.7998
$ws.Range("N134").Value = -16317.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5031.6
$ws.Range("J31").Value = 7242
$ws.Range("L31").Value = 7242
$ws.Range("N31").Value = -7832
$ws.Range("H34").Value = 5031.6
$ws.Range("J34").Value = 7242
$ws.Range("L34").Value = 7242
$ws.Range("N34").Value = -7646
$ws.Range("H58").Value = 3646.3914
$ws.Range("I58").Value = 1880.125
$ws.Range("J58").Value = 4588.4
$ws.Range("K58").Value = 1880.125
$ws.Range("L58").Value = 4588.4
$ws.Range("M58").Value = -1677.125
$ws.Range("N58").Value = -4994.4
$ws.Range("H105").Value = 1401.75
$ws.Range("I105").Value = 1401.75
$ws.Range("K105").Value = 1401.75
$ws.Range("M105").Value = 345.25
$ws.Range("H107").Value = 648.4
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()  # cached value removed in the source edit
$ws.Range("H132").Value = 2493.7778
$ws.Range("I132").Value = 2392.1667
$ws.Range("J132").Value = 2697
$ws.Range("K132").Value = 7176.500100000001
$ws.Range("L132").Value = 8091
$ws.Range("M132").Value = -4646.500100000001
$ws.Range("N132").Value = -13151
$ws.Range("H134").Value = 2438.2903
$ws.Range("I134").Value = 1848.174
$ws.Range("J134").Value = 4134.875
$ws.Range("K134").Value = 5544.522
$ws.Range("L134").Value = 12404.625
$ws.Range("M134").Value = -3009.522
$ws.Range("N134").Value = -17474.625
$ws.Range("H136").Value = 3646.3914
$ws.Range("I136").Value = 1880.125
$ws.Range("J136").Value = 4588.4
$ws.Range("K136").Value = 5640.375
$ws.Range("L136").Value = 13765.2
$ws.Range("M136").Value = -3090.375
$ws.Range("N136").Value = -18865.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 84515180
$ws.Range("I4").Value = 92198296
$ws.Range("K4").Value = 276594888
$ws.Range("M4").Value = -276594776
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7346.1816
$ws.Range("I70").Value = 5808
$ws.Range("J70").Value = 7500
$ws.Range("K70").Value = 5808
$ws.Range("L70").Value = 7500
$ws.Range("M70").Value = -5538
$ws.Range("N70").Value = -8040
$ws.Range("H73").Value = 7346.1816
$ws.Range("I73").Value = 5808
$ws.Range("J73").Value = 7500
$ws.Range("K73").Value = 5808
$ws.Range("L73").Value = 7500
$ws.Range("M73").Value = -4872
$ws.Range("N73").Value = -9372
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11379.6
$ws.Range("I22").Value = 1633
$ws.Range("J22").Value = 25999.5
$ws.Range("K22").Value = 1633
$ws.Range("L22").Value = 25999.5
$ws.Range("M22").Value = -1338
$ws.Range("N22").Value = -26589.5
$ws.Range("H27").Value = 11379.6
$ws.Range("I27").Value = 1633
$ws.Range("J27").Value = 25999.5
$ws.Range("K27").Value = 1633
$ws.Range("L27").Value = 25999.5
$ws.Range("M27").Value = -1526
$ws.Range("N27").Value = -26213.5
$ws.Range("H82").Value = 2817.862
$ws.Range("I82").Value = 3074.238
$ws.Range("J82").Value = 2144.875
$ws.Range("K82").Value = 3074.238
$ws.Range("L82").Value = 2144.875
$ws.Range("M82").Value = -2713.238
$ws.Range("N82").Value = -2866.875
$ws.Range("H85").Value = 2817.862
$ws.Range("I85").Value = 3074.238
$ws.Range("J85").Value = 2144.875
$ws.Range("K85").Value = 3074.238
$ws.Range("L85").Value = 2144.875
$ws.Range("M85").Value = -1826.238
$ws.Range("N85").Value = -4640.875
$ws.Range("H93").Value = 977.0909
$ws.Range("I93").Value = 976
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 976
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 272
$ws.Range("N93").Value = -3496
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2263.5
$ws.Range("I81").Value = 2263.5
$ws.Range("K81").Value = 4527
$ws.Range("M81").Value = -3466
$ws.Range("H84").Value = 2263.5
$ws.Range("I84").Value = 2263.5
$ws.Range("K84").Value = 22635
$ws.Range("M84").Value = -17331
$ws.Range("H122").Value = 3097.8
$ws.Range("I122").Value = 2997.25
$ws.Range("K122").Value = 8991.75
$ws.Range("M122").Value = -6541.75
$ws.Range("H132").Value = 1556.6
$ws.Range("I132").Value = 1537.2858
$ws.Range("K132").Value = 4611.857400000001
$ws.Range("M132").Value = -2081.857400000001
